# Reorder/drop columns in the sheet:
#  - drop old columns A (distance_km), B (match_score), C (match_rank)
#  - keep the remaining columns, but reorder the language columns
# New header order (A1:W1):
#   cand_gender, cand_age_bucket, cand_domicile_province, cand_domicile_region,
#   job_contract_type, job_work_province, Svedese, Spagnolo, Finlandese, Ebraico,
#   Persiano, Portoghese, Cinese, Arabo, Serbo, Albanese, Croato, Ceco, Danese,
#   Rumeno, Macedone, Tedesco, same_location

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Number of data rows (excluding header) and old/new column counts
$numDataRows = 10
$oldCols = 26
$newCols = 23
$numRows = 11

# 1-indexed old column number that each new column comes from
$mapping = @(4, 5, 6, 7, 8, 9, 14, 16, 18, 23, 13, 22, 25, 19, 21, 15, 17, 10, 20, 24, 11, 12, 26)

# Read all existing values (header + data) before mutating anything
$oldRange = $ws.Range("A1").Resize($numRows, $oldCols)
$oldVals = $oldRange.Value()

# Stash a copy of the header cell format (bold, centered, bordered) in a
# scratch cell far outside the used range, since Clear() below will wipe
# the formatting (and clipboard contents get invalidated by Clear too).
$scratch = $ws.Range("AZ1")
$ws.Range("A1").Copy()
$scratch.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Build the new grid in memory
$newVals = New-Object 'object[,]' $numRows, $newCols

for ($c = 1; $c -le $newCols; $c++) {
    $srcCol = $mapping[$c - 1]
    for ($r = 1; $r -le $numRows; $r++) {
        $newVals[$r - 1, $c - 1] = $oldVals[$r, $srcCol]
    }
}

# Clear the whole old used area, then write the new grid
$ws.Range("A1").Resize($numRows, $oldCols).Clear()

$destRange = $ws.Range("A1").Resize($numRows, $newCols)
$destRange.Value = $newVals

# Reapply the header style (bold, centered, bordered) that Clear() removed
$scratch.Copy()
$headerRow = $ws.Range("A1").Resize(1, $newCols)
$headerRow.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Remove the temporary scratch formatting
$scratch.Clear()
